# Apply cell value updates per the diff (crypto price/volume refresh, plus row 46
# insertion of BabyDogeCoin which shifts RenderToken/Aave/Cronos/Mantle/USDD down
# by one row and drops the old EnergySwap row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.240.13"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "1.608.38"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("E11").Value = "  +4.11%  "
$ws.Range("D12").Value = "1.832.51"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "1.608.30"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "26.224.52"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "198.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.130"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +1.35%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "1.108.87"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.501"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.780"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.15%  "
$ws.Range("D43").Value = "1.745.75"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0108"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.41%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.48%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.410"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "
